$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.506.83'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.571.53'
$ws.Range("E3").Value = '  -0.69%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '212.18'
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.494'
$ws.Range("E6").Value = '  -0.55%  '
$ws.Range("E7").Value = '  +0.31%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '46.27'
$ws.Range("E8").Value = '  +6.17%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '23.98'
$ws.Range("E9").Value = '  +0.88%  '
$ws.Range("E10").Value = '  -1.62%  '
$ws.Range("E11").Value = '  -1.48%  '
$ws.Range("E12").Value = '  -0.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.795.42'
$ws.Range("E13").Value = '  -0.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.568.29'
$ws.Range("E14").Value = '  -1.01%  '
$ws.Range("E15").Value = '  -1.67%  '
$ws.Range("E16").Value = '  -2.24%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '28.494.07'
$ws.Range("E17").Value = '  +0.53%  '
$ws.Range("E18").Value = '  -3.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '227.05'
$ws.Range("E19").Value = '  -2.48%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.35'
$ws.Range("E20").Value = '  -1.59%  '
$ws.Range("E22").Value = '  +0.23%  '
$ws.Range("E23").Value = '  -5.94%  '
$ws.Range("E24").Value = '  -2.56%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.08'
$ws.Range("E25").Value = '  +7.29%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '150.36'
$ws.Range("E26").Value = '  -0.73%  '
$ws.Range("E27").Value = '  -2.03%  '
$ws.Range("E28").Value = '  -2.69%  '
$ws.Range("E29").Value = '  -2.32%  '
$ws.Range("E30").Value = '  +0.16%  '
$ws.Range("E31").Value = '  -3.13%  '
$ws.Range("E32").Value = '  -1.58%  '
$ws.Range("E33").Value = '  -0.80%  '
$ws.Range("E34").Value = '  -0.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.397.77'
$ws.Range("E35").Value = '  -1.29%  '
$ws.Range("B36").Value = 'TrustWalletToken'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.04'
$ws.Range("E36").Value = '  -0.98%  '
$ws.Range("B37").Value = 'LidoDAOToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.55'
$ws.Range("E37").Value = '  -3.07%  '
$ws.Range("E38").Value = '  +1.87%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.60'
$ws.Range("E39").Value = '  +2.76%  '
$ws.Range("E40").Value = '  -0.70%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.534'
$ws.Range("E41").Value = '  -1.73%  '
$ws.Range("E42").Value = '  +0.26%  '
$ws.Range("E43").Value = '  -2.58%  '
$ws.Range("E44").Value = '  -1.25%  '
$ws.Range("E45").Value = '  +1.97%  '
$ws.Range("E46").Value = '  +0.77%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '62.97'
$ws.Range("E47").Value = '  -2.09%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.708.25'
$ws.Range("E48").Value = '  -0.53%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '86.09'
$ws.Range("E49").Value = '  -1.29%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0₆0103'
$ws.Range("E50").Value = '  -0.82%  '
$ws.Range("E51").Value = '  -1.43%  '
